# Append three more rows (12-14) to the "FirstSheet" table, repeating the
# same Abhi_0 / Abhi_1 / Abhi_2 pattern already used in rows 7-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 12; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = "Abhi_0"
    $ws.Cells.Item($row, 2).Value = "Abhi_1"
    $ws.Cells.Item($row, 3).Value = "Abhi_2"
}
